$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 115
$ws.Range("I61").Value = 115
$ws.Range("K61").Value = 345
$ws.Range("M61").Value = -173
$ws.Range("H94").Value = 428
$ws.Range("I94").Value = 428
$ws.Range("K94").Value = 428
$ws.Range("M94").Value = 23
$ws.Range("H135").Value = 1695.96
$ws.Range("I135").Value = 1059.4117
$ws.Range("J135").Value = 3048.625
$ws.Range("K135").Value = 9534.705300000001
$ws.Range("L135").Value = 27437.625
$ws.Range("M135").Value = -6999.705300000001
$ws.Range("N135").Value = -32507.625
$ws.Range("H138").Value = 6581847
$ws.Range("I138").Value = 1313.381
$ws.Range("K138").Value = 3940.143
$ws.Range("M138").Value = 1199.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 6033.857
$ws.Range("J16").Value = 7666
$ws.Range("L16").Value = 7666
$ws.Range("N16").Value = -8240
$ws.Range("H32").Value = 11828.805
$ws.Range("I32").Value = 6222.0938
$ws.Range("K32").Value = 6222.0938
$ws.Range("M32").Value = -5935.0938
$ws.Range("H110").Value = 6531.905
$ws.Range("I110").Value = 7253.8887
$ws.Range("K110").Value = 7253.8887
$ws.Range("M110").Value = -5208.8887
$ws.Range("H122").Value = 1464.8889
$ws.Range("I122").Value = 1310.8334
$ws.Range("J122").Value = 1773
$ws.Range("K122").Value = 3932.5002
$ws.Range("L122").Value = 5319
$ws.Range("M122").Value = -1482.5002
$ws.Range("N122").Value = -10219

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2268.2693
$ws.Range("I134").Value = 2229.02
$ws.Range("K134").Value = 6687.059999999999
$ws.Range("M134").Value = -4152.059999999999
$ws.Range("H135").Value = 58046.668
$ws.Range("J135").Value = 58046.668
$ws.Range("L135").Value = 58046.668
$ws.Range("N135").Value = -68186.668
$ws.Range("H138").Value = 62814.668
$ws.Range("J138").Value = 62814.668
$ws.Range("L138").Value = 62814.668
$ws.Range("N138").Value = -73094.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 37313.625
$ws.Range("I33").Value = 9265
$ws.Range("J33").Value = 46663.168
$ws.Range("K33").Value = 9265
$ws.Range("L33").Value = 46663.168
$ws.Range("M33").Value = -8886
$ws.Range("N33").Value = -47421.168
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = $null
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = $null
$ws.Range("H44").Value = 8466
$ws.Range("I44").Value = 8988
$ws.Range("J44").Value = 6900
$ws.Range("K44").Value = 8988
$ws.Range("L44").Value = 6900
$ws.Range("M44").Value = -8546
$ws.Range("N44").Value = -7784
$ws.Range("H86").Value = 5992.5713
$ws.Range("I86").Value = 5582.643
$ws.Range("J86").Value = 6812.4287
$ws.Range("K86").Value = 5582.643
$ws.Range("L86").Value = 6812.4287
$ws.Range("M86").Value = -4459.643
$ws.Range("N86").Value = -9058.4287
$ws.Range("H89").Value = 5992.5713
$ws.Range("I89").Value = 5582.643
$ws.Range("J89").Value = 6812.4287
$ws.Range("K89").Value = 27913.215
$ws.Range("L89").Value = 34062.14350000001
$ws.Range("M89").Value = -22297.215
$ws.Range("N89").Value = -45294.14350000001
$ws.Range("H124").Value = 44997
$ws.Range("J124").Value = 44997
$ws.Range("L124").Value = 44997
$ws.Range("N124").Value = -49907

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2276
$ws.Range("I114").Value = 1423.3334
$ws.Range("J114").Value = 3299.2
$ws.Range("K114").Value = 4270.0002
$ws.Range("L114").Value = 9897.599999999999
$ws.Range("M114").Value = -1016.0002
$ws.Range("N114").Value = -16405.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 37072036
$ws.Range("I18").Value = 55558056
$ws.Range("J18").Value = 100000
$ws.Range("K18").Value = 55558056
$ws.Range("L18").Value = 100000
$ws.Range("M18").Value = -55557763
$ws.Range("N18").Value = -100586
$ws.Range("H24").Value = 20000
$ws.Range("J24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("N24").Value = -20346
$ws.Range("H126").Value = 10003.241
$ws.Range("J126").Value = 3346.75
$ws.Range("L126").Value = 10040.25
$ws.Range("N126").Value = -14980.25
$ws.Range("H132").Value = 3402.3215
$ws.Range("I132").Value = 2933.8572
$ws.Range("K132").Value = 8801.5716
$ws.Range("M132").Value = -6271.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1557.4286
$ws.Range("I22").Value = 1135.8
$ws.Range("J22").Value = 1791.6666
$ws.Range("K22").Value = 1135.8
$ws.Range("L22").Value = 1791.6666
$ws.Range("M22").Value = -840.8
$ws.Range("N22").Value = -2381.6666
$ws.Range("H23").Value = 16720000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 16720000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 16720000
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = -16720460
$ws.Range("H27").Value = 1557.4286
$ws.Range("I27").Value = 1135.8
$ws.Range("J27").Value = 1791.6666
$ws.Range("K27").Value = 1135.8
$ws.Range("L27").Value = 1791.6666
$ws.Range("M27").Value = -1028.8
$ws.Range("N27").Value = -2005.6666
$ws.Range("H30").Value = 800000
$ws.Range("I30").Value = 800000
$ws.Range("K30").Value = 800000
$ws.Range("M30").Value = -799892
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("H46").Value = 3350
$ws.Range("I46").Value = 2700
$ws.Range("J46").Value = 4000
$ws.Range("K46").Value = 2700
$ws.Range("L46").Value = 4000
$ws.Range("M46").Value = -2512
$ws.Range("N46").Value = -4376
$ws.Range("H58").Value = 7985
$ws.Range("I58").Value = 7309.3335
$ws.Range("J58").Value = 8998.5
$ws.Range("K58").Value = 7309.3335
$ws.Range("L58").Value = 8998.5
$ws.Range("M58").Value = -7049.3335
$ws.Range("N58").Value = -9518.5
$ws.Range("H138").Value = 47430
$ws.Range("J138").Value = 45950
$ws.Range("L138").Value = 45950
$ws.Range("N138").Value = -56230

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3600
$ws.Range("I4").Value = 1800
$ws.Range("J4").Value = 4500
$ws.Range("K4").Value = 1800
$ws.Range("L4").Value = 4500
$ws.Range("M4").Value = -1687
$ws.Range("N4").Value = -4726
$ws.Range("H31").Value = 18000
$ws.Range("J31").Value = 18000
$ws.Range("L31").Value = 18000
$ws.Range("N31").Value = -18696
$ws.Range("H126").Value = 2262.2354
$ws.Range("I126").Value = 2187.4546
$ws.Range("K126").Value = 6562.3638
$ws.Range("M126").Value = -4092.3638
